# Add files via upload
# Fill in the next set of transactions (rows 10-17) on Sheet1, mirroring the
# pattern already present in rows 2-9, and move the active selection to D13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data to add: Date (serial), Name, Amount
$data = @(
    @(45323, "Gehalt", 3000),
    @(45324, "Miete", -800),
    @(45325, "Rewe", -100),
    @(45326, "Kino", -30),
    @(45327, "Rewe", -250),
    @(45328, "Werkstatt", -250),
    @(45329, "Restaurant", -100),
    @(45330, "Rewe", -80)
)

$startRow = 10
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $entry = $data[$i]

    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
}

$ws.Range("D13").Select() | Out-Null

$wb.Save()
